# Refresh cryptos list - Fri Aug 25 23:38:13 UTC 2023 (GitHub Actions).
# Source diff only touches columns D (Price) / E (Volume(1h)) for almost
# every row, plus a rank swap between BabyDogeCoin and Frax (rows 46/47).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Per-row Price (D) / Volume(1h) (E) updates ---
$ws.Range("D2").Value = "'26.097.53"
$ws.Range("E2").Value = "'  -0.58%  "
$ws.Range("D3").Value = "'1.657.06"
$ws.Range("E3").Value = "'  -0.34%  "
$ws.Range("E4").Value = "'  -0.39%  "
$ws.Range("D5").Value = "'218.60"
$ws.Range("E5").Value = "'  -0.25%  "
$ws.Range("D6").Value = "'0.5304"
$ws.Range("E6").Value = "'  +1.51%  "
$ws.Range("E7").Value = "'  -0.33%  "
$ws.Range("D8").Value = "'0.2618"
$ws.Range("E8").Value = "'  -2.04%  "
$ws.Range("D9").Value = "'0.06339"
$ws.Range("E9").Value = "'  -0.18%  "
$ws.Range("D10").Value = "'20.44"
$ws.Range("E10").Value = "'  -3.04%  "
$ws.Range("E11").Value = "'  +0.43%  "
$ws.Range("E12").Value = "'  +1.31%  "
$ws.Range("D13").Value = "'1.635.76"
$ws.Range("E13").Value = "'  -1.33%  "
$ws.Range("E14").Value = "'  +0.05%  "
$ws.Range("D15").Value = "'0.0₅8158"
$ws.Range("E15").Value = "'  -0.76%  "
$ws.Range("D16").Value = "'65.22"
$ws.Range("E16").Value = "'  +0.27%  "
$ws.Range("D17").Value = "'26.108.58"
$ws.Range("E17").Value = "'  -0.60%  "
$ws.Range("E18").Value = "'  -0.36%  "
$ws.Range("D19").Value = "'4.553"
$ws.Range("E19").Value = "'  -2.23%  "
$ws.Range("D20").Value = "'193.18"
$ws.Range("E20").Value = "'  -1.17%  "
$ws.Range("E21").Value = "'  -0.88%  "
$ws.Range("D22").Value = "'6.028"
$ws.Range("E22").Value = "'  -1.06%  "
$ws.Range("D23").Value = "'1.003"
$ws.Range("E23").Value = "'  -0.49%  "
$ws.Range("D24").Value = "'140.12"
$ws.Range("E24").Value = "'  +0.67%  "
$ws.Range("D25").Value = "'0.1246"
$ws.Range("E25").Value = "'  +0.02%  "
$ws.Range("D26").Value = "'7.273"
$ws.Range("E26").Value = "'  +0.53%  "
$ws.Range("D27").Value = "'16.17"
$ws.Range("E27").Value = "'  -0.16%  "
$ws.Range("E28").Value = "'  +1.13%  "
$ws.Range("D29").Value = "'0.05946"
$ws.Range("E29").Value = "'  -0.54%  "
$ws.Range("D30").Value = "'1.276"
$ws.Range("E30").Value = "'  -0.62%  "
$ws.Range("D31").Value = "'3.510"
$ws.Range("E31").Value = "'  -2.90%  "
$ws.Range("D32").Value = "'3.238"
$ws.Range("E32").Value = "'  -1.60%  "
$ws.Range("D33").Value = "'1.560"
$ws.Range("E33").Value = "'  -4.44%  "
$ws.Range("D34").Value = "'0.9501"
$ws.Range("E34").Value = "'  -3.27%  "
$ws.Range("E35").Value = "'  -0.61%  "
$ws.Range("D36").Value = "'2.769"
$ws.Range("E36").Value = "'  -0.53%  "
$ws.Range("D37").Value = "'0.5641"
$ws.Range("E37").Value = "'  -4.40%  "
$ws.Range("D38").Value = "'0.01611"
$ws.Range("E38").Value = "'  +0.74%  "
$ws.Range("D39").Value = "'5.835"
$ws.Range("E39").Value = "'  -2.92%  "
$ws.Range("D40").Value = "'0.8440"
$ws.Range("E40").Value = "'  -1.88%  "
$ws.Range("E41").Value = "'  -0.30%  "
$ws.Range("D42").Value = "'101.47"
$ws.Range("E42").Value = "'  +1.38%  "
$ws.Range("D43").Value = "'1.013.20"
$ws.Range("E43").Value = "'  -1.70%  "
$ws.Range("D44").Value = "'1.801.26"
$ws.Range("E44").Value = "'  -0.22%  "
$ws.Range("D45").Value = "'57.09"
$ws.Range("E45").Value = "'  -0.71%  "
$ws.Range("E48").Value = "'  +1.28%  "
$ws.Range("D49").Value = "'0.05155"
$ws.Range("E49").Value = "'  -0.58%  "
$ws.Range("D50").Value = "'1.470"
$ws.Range("E50").Value = "'  -0.17%  "
$ws.Range("D51").Value = "'7.725"
$ws.Range("E51").Value = "'  -4.32%  "

# --- Rows 46/47: BabyDogeCoin overtakes Frax in rank order ---
$ws.Range("B46").Value = "'BabyDogeCoin"
$ws.Range("C46").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.0₈104"
$ws.Range("E46").Value = "'  -7.09%  "

$ws.Range("B47").Value = "'Frax"
$ws.Range("C47").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.001"
$ws.Range("E47").Value = "'  -1.04%  "

